$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a brand-new "Player Info" sheet in front of everything
#    else (Worksheets.Add() with no args inserts before the active
#    sheet, which is the first / left-most tab).
# ------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add()
$infoSheet.Name = "Player Info"

# ------------------------------------------------------------------
# 2) Re-fetch the pre-existing sheets by name now that the tab order
#    has shifted (sheet references are resolved by position, so they
#    must be looked up again after the insert above).
# ------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# Header row
$infoSheet.Cells.Item(1, 1).Value = "ID"
$infoSheet.Cells.Item(1, 2).Value = "NAME"
$infoSheet.Cells.Item(1, 3).Value = "BATTING_HAND"
$infoSheet.Cells.Item(1, 4).Value = "BOWL_STYLE"

# Match the bold / centered / bordered look used for every other
# header row in this workbook.
$hdr = $infoSheet.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4160    # xlTop
$hdr.Borders.LineStyle = 1

# Data row
$infoSheet.Cells.Item(2, 2).Value = "Matthew James Humphreys"
$infoSheet.Cells.Item(2, 3).Value = "Right Handed"
$infoSheet.Cells.Item(2, 4).Value = "Left Arm Orthodox"

# ID is numeric-looking text ("7189"), so force a text format before
# assigning it, otherwise it gets auto-coerced into a number.
$idCell = $infoSheet.Cells.Item(2, 1)
$idCell.NumberFormat = "@"
$idCell.Value = "7189"

[void]$infoSheet.Range("A1").Select()

# ------------------------------------------------------------------
# 3) Rename the MATCH_CARD_LINK column to MATCH_CODE on both the
#    batting and bowling sheets, and replace the full scorecard URL
#    with just the numeric match code.
# ------------------------------------------------------------------

# --- ODI Batting: MATCH_CARD_LINK lives in column D ---
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$batD2 = $battingSheet.Cells.Item(2, 4)
$batD2.NumberFormat = "@"
$batD2.Value = "4729"

$batD3 = $battingSheet.Cells.Item(3, 4)
$batD3.NumberFormat = "@"
$batD3.Value = "4734"

# --- ODI Bowling: MATCH_CARD_LINK lives in column B ---
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowB2 = $bowlingSheet.Cells.Item(2, 2)
$bowB2.NumberFormat = "@"
$bowB2.Value = "4729"

$bowB3 = $bowlingSheet.Cells.Item(3, 2)
$bowB3.NumberFormat = "@"
$bowB3.Value = "4734"

Write-Host "Done. Sheets:" ($wb.Worksheets | ForEach-Object { $_.Name })
